$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "HOLI50"
$ws.Range("D2").Value = "8/2/2025, 11:37:42 am"

# Row 3
$ws.Range("B3").Value = "HOLI50"
$ws.Range("D3").Value = "8/2/2025, 11:38:21 am"

# Row 4
$ws.Range("B4").Value = "HOLI50"
$ws.Range("C4").Value = "Karelibagh Branch "
$ws.Range("D4").Value = "8/2/2025, 3:02:38 pm"
$ws.Range("E4").Value = "Dev"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "9898561251"

# Row 5
$ws.Range("B5").Value = "HOLI50"
$ws.Range("C5").Value = "Karelibagh Branch "
$ws.Range("D5").Value = "8/2/2025, 3:03:07 pm"
$ws.Range("E5").Value = "Poojan"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "8160754098"
